$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, derived from the updated symbol-list scrape.
$updates = @{
    2 = @{ "D" = "260.82"; "E" = "-0.34%" }
    3 = @{ "D" = "27.53"; "E" = "0.12%" }
    4 = @{ "D" = "4.715"; "E" = "-0.84%" }
    5 = @{ "D" = "0.06219"; "E" = "2.43%" }
    6 = @{ "D" = "6.734"; "E" = "0.37%" }
    7 = @{ "D" = "0.8498"; "E" = "-1.47%" }
    8 = @{ "D" = "0.9095"; "E" = "-1.48%" }
    9 = @{ "D" = "0.1399"; "E" = "-0.75%" }
    10 = @{ "D" = "0.04816"; "E" = "-6.10%" }
    11 = @{ "D" = "0.07083"; "E" = "-0.93%" }
    12 = @{ "D" = "0.03123"; "E" = "1.85%" }
    13 = @{ "D" = "0.09056"; "E" = "-0.65%" }
    14 = @{ "D" = "0.001525"; "E" = "-0.76%" }
    15 = @{ "D" = "0.0006143"; "E" = "0.98%" }
    16 = @{ "D" = "0.006059"; "E" = "-2.16%" }
    17 = @{ "E" = "0.05%" }
    18 = @{ "D" = "3.169"; "E" = "0.12%" }
    19 = @{ "E" = "-0.51%" }
    21 = @{ "E" = "1.60%" }
    22 = @{ "D" = "4.091"; "E" = "-0.06%" }
    23 = @{ "D" = "0.04260"; "E" = "0.16%" }
    24 = @{ "D" = "0.001222"; "E" = "0.32%" }
    25 = @{ "D" = "0.004083"; "E" = "4.33%" }
    26 = @{ "E" = "0.05%" }
    27 = @{ "E" = "4.40%" }
    40 = @{ "D" = "0.03904"; "E" = "0.64%" }
    41 = @{ "D" = "0.1110"; "E" = "-0.52%" }
    42 = @{ "D" = "0.004113"; "E" = "-0.44%" }
    43 = @{ "D" = "0.002147"; "E" = "-2.78%" }
    44 = @{ "D" = "0.01344"; "E" = "-11.61%" }
    45 = @{ "D" = "0.00005130"; "E" = "-2.89%" }
    46 = @{ "D" = "0.00000000750"; "E" = "0.03%" }
    47 = @{ "D" = "0.03402"; "E" = "-37.66%" }
    48 = @{ "D" = "0.06503"; "E" = "-50.78%" }
    49 = @{ "D" = "0.00002101"; "E" = "0.03%" }
    50 = @{ "D" = "0.0002001"; "E" = "0.03%" }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cellAddr = "$col$row"
        $cell = $ws.Range($cellAddr)
        # Force text storage so numeric-looking strings (prices, percentages)
        # are kept verbatim instead of being reinterpreted as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$row][$col]
        $cell.Style = "Normal"
    }
}
